$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 69
$ws.Range("H69").Value = 3780
$ws.Range("J69").Value = 3780
$ws.Range("L69").Value = 11340
$ws.Range("N69").Value = -13088

# Row 72
$ws.Range("H72").Value = 3780
$ws.Range("J72").Value = 3780
$ws.Range("L72").Value = 34020
$ws.Range("N72").Value = -42756

# Row 96
$ws.Range("H96").Value = 597.4286
$ws.Range("I96").Value = 597.4286
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1792.2858
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -419.2857999999999

# Row 99
$ws.Range("H99").Value = 1357.7142
$ws.Range("I99").Value = 885.6
$ws.Range("J99").Value = 1620
$ws.Range("K99").Value = 2656.8
$ws.Range("L99").Value = 4860
$ws.Range("M99").Value = -1158.8
$ws.Range("N99").Value = -7856

# Row 119
$ws.Range("H119").Value = 652.5
$ws.Range("J119").Value = 652.5
$ws.Range("L119").Value = 1957.5
$ws.Range("N119").Value = -11633.5

# Row 129
$ws.Range("H129").Value = 947.50793
$ws.Range("I129").Value = 403.42856
$ws.Range("J129").Value = 1102.9592
$ws.Range("K129").Value = 1210.28568
$ws.Range("L129").Value = 3308.8776
$ws.Range("M129").Value = 3789.71432
$ws.Range("N129").Value = -13308.8776

# Row 132
$ws.Range("H132").Value = 1280.1519
$ws.Range("I132").Value = 1112.7606
$ws.Range("J132").Value = 2765.75
$ws.Range("K132").Value = 3338.2818
$ws.Range("L132").Value = 8297.25
$ws.Range("M132").Value = -808.2818000000002
$ws.Range("N132").Value = -13357.25

# Row 135
$ws.Range("H135").Value = 1049.8
$ws.Range("I135").Value = 906.3043
$ws.Range("J135").Value = 2700
$ws.Range("K135").Value = 8156.7387
$ws.Range("L135").Value = 24300
$ws.Range("M135").Value = -5621.7387
$ws.Range("N135").Value = -29370

# Row 138
$ws.Range("H138").Value = 2044622.9
$ws.Range("I138").Value = 5265335
$ws.Range("J138").Value = 4838.4
$ws.Range("K138").Value = 15796005
$ws.Range("L138").Value = 14515.2
$ws.Range("M138").Value = -15790865
$ws.Range("N138").Value = -24795.2

# Row 141
$ws.Range("H141").Value = 9668.125
$ws.Range("I141").Value = 5281.8423
$ws.Range("K141").Value = 15845.5269
$ws.Range("M141").Value = -10665.5269

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 10342.4
$ws.Range("I63").Value = 11996.25
$ws.Range("J63").Value = 9239.833000000001
$ws.Range("K63").Value = 11996.25
$ws.Range("L63").Value = 9239.833000000001
$ws.Range("M63").Value = -11310.25
$ws.Range("N63").Value = -10611.833

# Row 66
$ws.Range("H66").Value = 10342.4
$ws.Range("I66").Value = 11996.25
$ws.Range("J66").Value = 9239.833000000001
$ws.Range("K66").Value = 59981.25
$ws.Range("L66").Value = 46199.165
$ws.Range("M66").Value = -56549.25
$ws.Range("N66").Value = -53063.165

# Row 110
$ws.Range("H110").Value = 1989.6666
$ws.Range("I110").Value = 1989.6666
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1989.6666
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").Value = 55.33339999999998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 31022.47
$ws.Range("I86").Value = 1434.1111
$ws.Range("J86").Value = 64309.375
$ws.Range("K86").Value = 1434.1111
$ws.Range("L86").Value = 64309.375
$ws.Range("M86").Value = -311.1111000000001
$ws.Range("N86").Value = -66555.375

# Row 89
$ws.Range("H89").Value = 31022.47
$ws.Range("I89").Value = 1434.1111
$ws.Range("J89").Value = 64309.375
$ws.Range("K89").Value = 7170.5555
$ws.Range("L89").Value = 321546.875
$ws.Range("M89").Value = -1554.5555
$ws.Range("N89").Value = -332778.875

# Row 132
$ws.Range("H132").Value = 71340
$ws.Range("J132").Value = 71340
$ws.Range("L132").Value = 71340
$ws.Range("N132").Value = -81460

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 86
$ws.Range("H86").Value = 1696
$ws.Range("I86").Value = 1335.6666
$ws.Range("J86").Value = 2777
$ws.Range("K86").Value = 1335.6666
$ws.Range("L86").Value = 2777
$ws.Range("M86").Value = -212.6666
$ws.Range("N86").Value = -5023

# Row 89
$ws.Range("H89").Value = 1696
$ws.Range("I89").Value = 1335.6666
$ws.Range("J89").Value = 2777
$ws.Range("K89").Value = 6678.333000000001
$ws.Range("L89").Value = 13885
$ws.Range("M89").Value = -1062.333000000001
$ws.Range("N89").Value = -25117

# Row 132
$ws.Range("H132").Value = 367226.7
$ws.Range("I132").Value = 542313.8
$ws.Range("J132").Value = 2461.8333
$ws.Range("K132").Value = 1626941.4
$ws.Range("L132").Value = 7385.499899999999
$ws.Range("M132").Value = -1624411.4
$ws.Range("N132").Value = -12445.4999

# Row 141
$ws.Range("H141").Value = 56500.5
$ws.Range("I141").Value = 33001
$ws.Range("J141").Value = 80000
$ws.Range("K141").Value = 33001
$ws.Range("L141").Value = 80000
$ws.Range("M141").Value = -27821
$ws.Range("N141").Value = -90360

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 44
$ws.Range("H44").Value = 660
$ws.Range("I44").Value = 300
$ws.Range("J44").Value = 732
$ws.Range("K44").Value = 900
$ws.Range("L44").Value = 2196
$ws.Range("M44").Value = -502
$ws.Range("N44").Value = -2992

# Row 141
$ws.Range("H141").Value = 3848.7856
$ws.Range("I141").Value = 3834.6155
$ws.Range("J141").Value = 4033
$ws.Range("K141").Value = 11503.8465
$ws.Range("L141").Value = 12099
$ws.Range("M141").Value = -6323.8465
$ws.Range("N141").Value = -22459

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6542.647
$ws.Range("I70").Value = 5343.2
$ws.Range("J70").Value = 7042.4165
$ws.Range("K70").Value = 5343.2
$ws.Range("L70").Value = 7042.4165
$ws.Range("M70").Value = -5073.2
$ws.Range("N70").Value = -7582.4165

# Row 73
$ws.Range("H73").Value = 6542.647
$ws.Range("I73").Value = 5343.2
$ws.Range("J73").Value = 7042.4165
$ws.Range("K73").Value = 5343.2
$ws.Range("L73").Value = 7042.4165
$ws.Range("M73").Value = -4407.2
$ws.Range("N73").Value = -8914.416499999999

# Row 80
$ws.Range("H80").Value = 3676.375
$ws.Range("I80").Value = 3001
$ws.Range("J80").Value = 4802
$ws.Range("K80").Value = 3001
$ws.Range("L80").Value = 4802
$ws.Range("M80").Value = -2003
$ws.Range("N80").Value = -6798

# Row 83
$ws.Range("H83").Value = 3676.375
$ws.Range("I83").Value = 3001
$ws.Range("J83").Value = 4802
$ws.Range("K83").Value = 15005
$ws.Range("L83").Value = 24010
$ws.Range("M83").Value = -10013
$ws.Range("N83").Value = -33994

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 8351.727999999999
$ws.Range("I132").Value = 8420.8125
$ws.Range("J132").Value = 8167.5
$ws.Range("K132").Value = 25262.4375
$ws.Range("L132").Value = 24502.5
$ws.Range("M132").Value = -22732.4375
$ws.Range("N132").Value = -29562.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 133
$ws.Range("H133").Value = 46150
$ws.Range("J133").Value = 46150
$ws.Range("L133").Value = 46150
$ws.Range("N133").Value = -56270
